# Update each arithmetic-problem cell in the practice-sheet table with its
# new expression. Every "old" value below is unique in the document, so a
# simple MatchCase + MatchWholeWord Find/Replace on $d.Content safely
# targets exactly one <w:t> run each, regardless of table/cell position.
$d = $word.ActiveDocument

$d.Content.Find.Execute("37+7=", $true, $true, $false, $false, $false, $true, 1, $false, "26+37=", 2) | Out-Null
$d.Content.Find.Execute("91-34=", $true, $true, $false, $false, $false, $true, 1, $false, "19+29=", 2) | Out-Null
$d.Content.Find.Execute("61-38=", $true, $true, $false, $false, $false, $true, 1, $false, "24+58=", 2) | Out-Null
$d.Content.Find.Execute("47+35=", $true, $true, $false, $false, $false, $true, 1, $false, "45+37=", 2) | Out-Null
$d.Content.Find.Execute("65-46=", $true, $true, $false, $false, $false, $true, 1, $false, "77+19=", 2) | Out-Null
$d.Content.Find.Execute("17+39=", $true, $true, $false, $false, $false, $true, 1, $false, "97-88=", 2) | Out-Null
$d.Content.Find.Execute("61-39=", $true, $true, $false, $false, $false, $true, 1, $false, "64-36=", 2) | Out-Null
$d.Content.Find.Execute("55+9=", $true, $true, $false, $false, $false, $true, 1, $false, "80-75=", 2) | Out-Null
$d.Content.Find.Execute("91-52=", $true, $true, $false, $false, $false, $true, 1, $false, "63+9=", 2) | Out-Null
$d.Content.Find.Execute("49+23=", $true, $true, $false, $false, $false, $true, 1, $false, "66+6=", 2) | Out-Null
$d.Content.Find.Execute("65-7=", $true, $true, $false, $false, $false, $true, 1, $false, "28+58=", 2) | Out-Null
$d.Content.Find.Execute("40-14=", $true, $true, $false, $false, $false, $true, 1, $false, "34+7=", 2) | Out-Null
$d.Content.Find.Execute("43-9=", $true, $true, $false, $false, $false, $true, 1, $false, "92-48=", 2) | Out-Null
$d.Content.Find.Execute("31-3=", $true, $true, $false, $false, $false, $true, 1, $false, "67+26=", 2) | Out-Null
$d.Content.Find.Execute("27+7=", $true, $true, $false, $false, $false, $true, 1, $false, "27+68=", 2) | Out-Null
$d.Content.Find.Execute("17+44=", $true, $true, $false, $false, $false, $true, 1, $false, "71-58=", 2) | Out-Null
$d.Content.Find.Execute("4+68=", $true, $true, $false, $false, $false, $true, 1, $false, "77+5=", 2) | Out-Null
$d.Content.Find.Execute("55-19=", $true, $true, $false, $false, $false, $true, 1, $false, "30-29=", 2) | Out-Null
$d.Content.Find.Execute("60-47=", $true, $true, $false, $false, $false, $true, 1, $false, "14+68=", 2) | Out-Null
$d.Content.Find.Execute("8+86=", $true, $true, $false, $false, $false, $true, 1, $false, "48+49=", 2) | Out-Null
$d.Content.Find.Execute("23+58=", $true, $true, $false, $false, $false, $true, 1, $false, "61-26=", 2) | Out-Null
$d.Content.Find.Execute("64-18=", $true, $true, $false, $false, $false, $true, 1, $false, "52+29=", 2) | Out-Null
$d.Content.Find.Execute("82-44=", $true, $true, $false, $false, $false, $true, 1, $false, "15+56=", 2) | Out-Null
$d.Content.Find.Execute("13-8=", $true, $true, $false, $false, $false, $true, 1, $false, "69+9=", 2) | Out-Null
$d.Content.Find.Execute("91-29=", $true, $true, $false, $false, $false, $true, 1, $false, "38+48=", 2) | Out-Null
$d.Content.Find.Execute("77-68=", $true, $true, $false, $false, $false, $true, 1, $false, "81-42=", 2) | Out-Null
$d.Content.Find.Execute("42-18=", $true, $true, $false, $false, $false, $true, 1, $false, "91-25=", 2) | Out-Null
$d.Content.Find.Execute("71-15=", $true, $true, $false, $false, $false, $true, 1, $false, "28+4=", 2) | Out-Null
$d.Content.Find.Execute("54-48=", $true, $true, $false, $false, $false, $true, 1, $false, "28+45=", 2) | Out-Null
$d.Content.Find.Execute("91-12=", $true, $true, $false, $false, $false, $true, 1, $false, "29+25=", 2) | Out-Null
$d.Content.Find.Execute("87-69=", $true, $true, $false, $false, $false, $true, 1, $false, "21-14=", 2) | Out-Null
$d.Content.Find.Execute("29+8=", $true, $true, $false, $false, $false, $true, 1, $false, "92-54=", 2) | Out-Null
$d.Content.Find.Execute("34+28=", $true, $true, $false, $false, $false, $true, 1, $false, "82-8=", 2) | Out-Null
$d.Content.Find.Execute("56+36=", $true, $true, $false, $false, $false, $true, 1, $false, "51-37=", 2) | Out-Null
$d.Content.Find.Execute("52-8=", $true, $true, $false, $false, $false, $true, 1, $false, "7+7=", 2) | Out-Null
$d.Content.Find.Execute("93-24=", $true, $true, $false, $false, $false, $true, 1, $false, "62-53=", 2) | Out-Null
$d.Content.Find.Execute("18+39=", $true, $true, $false, $false, $false, $true, 1, $false, "26+49=", 2) | Out-Null
$d.Content.Find.Execute("64+7=", $true, $true, $false, $false, $false, $true, 1, $false, "30-1=", 2) | Out-Null
$d.Content.Find.Execute("45+7=", $true, $true, $false, $false, $false, $true, 1, $false, "59+17=", 2) | Out-Null
$d.Content.Find.Execute("22+59=", $true, $true, $false, $false, $false, $true, 1, $false, "84-36=", 2) | Out-Null
$d.Content.Find.Execute("74-29=", $true, $true, $false, $false, $false, $true, 1, $false, "74-9=", 2) | Out-Null
$d.Content.Find.Execute("51-48=", $true, $true, $false, $false, $false, $true, 1, $false, "42-17=", 2) | Out-Null
$d.Content.Find.Execute("27+65=", $true, $true, $false, $false, $false, $true, 1, $false, "70-26=", 2) | Out-Null
$d.Content.Find.Execute("42+39=", $true, $true, $false, $false, $false, $true, 1, $false, "3+18=", 2) | Out-Null
$d.Content.Find.Execute("23+18=", $true, $true, $false, $false, $false, $true, 1, $false, "82-78=", 2) | Out-Null
$d.Content.Find.Execute("19+56=", $true, $true, $false, $false, $false, $true, 1, $false, "18+17=", 2) | Out-Null
$d.Content.Find.Execute("80-2=", $true, $true, $false, $false, $false, $true, 1, $false, "60-11=", 2) | Out-Null
$d.Content.Find.Execute("89+5=", $true, $true, $false, $false, $false, $true, 1, $false, "87-58=", 2) | Out-Null
$d.Content.Find.Execute("13-6=", $true, $true, $false, $false, $false, $true, 1, $false, "92-47=", 2) | Out-Null
$d.Content.Find.Execute("56-39=", $true, $true, $false, $false, $false, $true, 1, $false, "73-68=", 2) | Out-Null
$d.Content.Find.Execute("48+25=", $true, $true, $false, $false, $false, $true, 1, $false, "64-28=", 2) | Out-Null
$d.Content.Find.Execute("68+26=", $true, $true, $false, $false, $false, $true, 1, $false, "73-47=", 2) | Out-Null
$d.Content.Find.Execute("15+38=", $true, $true, $false, $false, $false, $true, 1, $false, "48-19=", 2) | Out-Null
$d.Content.Find.Execute("44-26=", $true, $true, $false, $false, $false, $true, 1, $false, "81-36=", 2) | Out-Null
$d.Content.Find.Execute("65-58=", $true, $true, $false, $false, $false, $true, 1, $false, "84-49=", 2) | Out-Null
$d.Content.Find.Execute("7+75=", $true, $true, $false, $false, $false, $true, 1, $false, "53-37=", 2) | Out-Null
$d.Content.Find.Execute("85-49=", $true, $true, $false, $false, $false, $true, 1, $false, "90-82=", 2) | Out-Null
$d.Content.Find.Execute("55-29=", $true, $true, $false, $false, $false, $true, 1, $false, "72-58=", 2) | Out-Null
$d.Content.Find.Execute("55+16=", $true, $true, $false, $false, $false, $true, 1, $false, "26+16=", 2) | Out-Null
$d.Content.Find.Execute("33-7=", $true, $true, $false, $false, $false, $true, 1, $false, "38+55=", 2) | Out-Null
$d.Content.Find.Execute("16+19=", $true, $true, $false, $false, $false, $true, 1, $false, "82-19=", 2) | Out-Null
$d.Content.Find.Execute("46-29=", $true, $true, $false, $false, $false, $true, 1, $false, "41-34=", 2) | Out-Null
$d.Content.Find.Execute("81-3=", $true, $true, $false, $false, $false, $true, 1, $false, "61-28=", 2) | Out-Null
$d.Content.Find.Execute("17+15=", $true, $true, $false, $false, $false, $true, 1, $false, "28+29=", 2) | Out-Null
$d.Content.Find.Execute("95-86=", $true, $true, $false, $false, $false, $true, 1, $false, "5+26=", 2) | Out-Null
$d.Content.Find.Execute("27+38=", $true, $true, $false, $false, $false, $true, 1, $false, "95-19=", 2) | Out-Null
$d.Content.Find.Execute("28+48=", $true, $true, $false, $false, $false, $true, 1, $false, "18+59=", 2) | Out-Null
$d.Content.Find.Execute("28+54=", $true, $true, $false, $false, $false, $true, 1, $false, "71-13=", 2) | Out-Null
$d.Content.Find.Execute("13+49=", $true, $true, $false, $false, $false, $true, 1, $false, "3+78=", 2) | Out-Null
$d.Content.Find.Execute("44-35=", $true, $true, $false, $false, $false, $true, 1, $false, "6+85=", 2) | Out-Null
$d.Content.Find.Execute("22+39=", $true, $true, $false, $false, $false, $true, 1, $false, "50-39=", 2) | Out-Null
$d.Content.Find.Execute("78-9=", $true, $true, $false, $false, $false, $true, 1, $false, "4+57=", 2) | Out-Null
$d.Content.Find.Execute("86+7=", $true, $true, $false, $false, $false, $true, 1, $false, "23-6=", 2) | Out-Null
$d.Content.Find.Execute("95-88=", $true, $true, $false, $false, $false, $true, 1, $false, "27+37=", 2) | Out-Null
$d.Content.Find.Execute("40-9=", $true, $true, $false, $false, $false, $true, 1, $false, "35+6=", 2) | Out-Null
$d.Content.Find.Execute("72-66=", $true, $true, $false, $false, $false, $true, 1, $false, "14-8=", 2) | Out-Null
$d.Content.Find.Execute("15+18=", $true, $true, $false, $false, $false, $true, 1, $false, "18+23=", 2) | Out-Null
$d.Content.Find.Execute("29+66=", $true, $true, $false, $false, $false, $true, 1, $false, "87+9=", 2) | Out-Null
$d.Content.Find.Execute("28+25=", $true, $true, $false, $false, $false, $true, 1, $false, "17+27=", 2) | Out-Null
$d.Content.Find.Execute("74-5=", $true, $true, $false, $false, $false, $true, 1, $false, "17+38=", 2) | Out-Null
$d.Content.Find.Execute("28+44=", $true, $true, $false, $false, $false, $true, 1, $false, "56+16=", 2) | Out-Null
$d.Content.Find.Execute("8+79=", $true, $true, $false, $false, $false, $true, 1, $false, "60-28=", 2) | Out-Null
$d.Content.Find.Execute("34+39=", $true, $true, $false, $false, $false, $true, 1, $false, "96-9=", 2) | Out-Null
$d.Content.Find.Execute("14+59=", $true, $true, $false, $false, $false, $true, 1, $false, "16+56=", 2) | Out-Null
$d.Content.Find.Execute("69+14=", $true, $true, $false, $false, $false, $true, 1, $false, "56+39=", 2) | Out-Null
$d.Content.Find.Execute("66-17=", $true, $true, $false, $false, $false, $true, 1, $false, "19+69=", 2) | Out-Null
$d.Content.Find.Execute("33-26=", $true, $true, $false, $false, $false, $true, 1, $false, "58+29=", 2) | Out-Null
$d.Content.Find.Execute("81-23=", $true, $true, $false, $false, $false, $true, 1, $false, "57-8=", 2) | Out-Null
$d.Content.Find.Execute("65+28=", $true, $true, $false, $false, $false, $true, 1, $false, "81-54=", 2) | Out-Null
$d.Content.Find.Execute("70-3=", $true, $true, $false, $false, $false, $true, 1, $false, "50-14=", 2) | Out-Null
$d.Content.Find.Execute("67+28=", $true, $true, $false, $false, $false, $true, 1, $false, "33-25=", 2) | Out-Null
$d.Content.Find.Execute("79+6=", $true, $true, $false, $false, $false, $true, 1, $false, "70-62=", 2) | Out-Null
$d.Content.Find.Execute("11-8=", $true, $true, $false, $false, $false, $true, 1, $false, "25+59=", 2) | Out-Null
$d.Content.Find.Execute("41-27=", $true, $true, $false, $false, $false, $true, 1, $false, "51-32=", 2) | Out-Null
$d.Content.Find.Execute("8+36=", $true, $true, $false, $false, $false, $true, 1, $false, "48+49=", 2) | Out-Null
$d.Content.Find.Execute("17+16=", $true, $true, $false, $false, $false, $true, 1, $false, "85-27=", 2) | Out-Null
$d.Content.Find.Execute("56-48=", $true, $true, $false, $false, $false, $true, 1, $false, "57+38=", 2) | Out-Null
$d.Content.Find.Execute("68+8=", $true, $true, $false, $false, $false, $true, 1, $false, "60-22=", 2) | Out-Null
$d.Content.Find.Execute("53-4=", $true, $true, $false, $false, $false, $true, 1, $false, "90-22=", 2) | Out-Null
$d.Content.Find.Execute("21-2=", $true, $true, $false, $false, $false, $true, 1, $false, "91-57=", 2) | Out-Null
